# Data/EC/NIT-9011377417.xlsx
# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# -> Replace the whole worker/mora table (rows 16-23) with a new, larger
#    table (rows 16-28: 7 workers x up to 3 periods each = 13 data rows),
#    update the summary totals, and push the two signature/footer rows
#    down from 28/29 to 33/34.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlPasteAll = -4104

# ---------------------------------------------------------------------
# 1) Relocate the footer (signature) block from rows 28-29 to rows 33-34,
#    preserving formatting+values cell by cell (copying merged ranges in
#    one shot confuses the style resolver, so do it column-by-column).
# ---------------------------------------------------------------------
foreach ($col in @("B","C","H","I","J")) {
    $ws.Range($col + "28").Copy()
    $ws.Range($col + "33").PasteSpecial($xlPasteFormats)

    $ws.Range($col + "29").Copy()
    $ws.Range($col + "34").PasteSpecial($xlPasteFormats)
}

$ws.Range("B33:C33").Merge()
$ws.Range("B34:C34").Merge()
$ws.Range("H33:J33").Merge()
$ws.Range("H34:J34").Merge()

$ws.Range("B33").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H33").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("B34").Value = "___________________________________"
$ws.Range("H34").Value = "FIRMA DEL REPRESENTANTE LEGAL"

# Clear out the old footer rows (28-29) - they become data rows / blank.
$ws.Range("B28:C28").UnMerge()
$ws.Range("H28:J28").UnMerge()
$ws.Range("B29:C29").UnMerge()
$ws.Range("H29:J29").UnMerge()
$ws.Range("B28:J29").Clear()

# ---------------------------------------------------------------------
# 2) Grow the data table from 8 rows (16-23) to 13 rows (16-28).
#    Row 23 currently carries the special "last row" (bottom border)
#    style - capture that into the new last row (28) first, then stamp
#    the regular "middle row" style (from row 22) across rows 23-27.
# ---------------------------------------------------------------------
$ws.Range("B23:J23").Copy()
$ws.Range("B28:J28").PasteSpecial($xlPasteFormats)

$ws.Range("B22:J22").Copy()
$ws.Range("B23:J27").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# 3) Write the new worker / mora-period data (rows 16-28).
#    Columns: B=Tipo Doc, C=N Doc, D=Nombre, E=Periodo, F=Valor Mora,
#    G=Salario Basico.
# ---------------------------------------------------------------------
$rows = @(
    @(16,"CC","8854791","ALEXANDER MARTINEZ ARIAS","2507",56940,1423500),
    @(17,"CC","9146396","RICARDO DE AVILA RAMOS","2507",56940,877803),
    @(18,"CC","33227517","DIANA MARCELA SALCEDO CASTELLAR","2507",56940,877803),
    @(19,"CC","1128054052","DAVID LEONARDO BARBOSA ESPITIA","2507",56940,1423500),
    @(20,"CC","1049829535","ANGEL DAVID ACEVEDO VICTOR","2507",56940,1423500),
    @(21,"CC","1049829535","ANGEL DAVID ACEVEDO VICTOR","2506",56940,1423500),
    @(22,"CC","1049829535","ANGEL DAVID ACEVEDO VICTOR","2505",45552,1423500),
    @(23,"CC","1002315450","YEIFER BARONA MANJARREZ","2507",56940,1423500),
    @(24,"CC","1002315450","YEIFER BARONA MANJARREZ","2506",56940,1423500),
    @(25,"CC","1002315450","YEIFER BARONA MANJARREZ","2505",45552,1423500),
    @(26,"CC","1051886922","LUIS MIGUEL MERCADO JUNCO","2507",56940,1423500),
    @(27,"CC","1051886922","LUIS MIGUEL MERCADO JUNCO","2506",56940,1423500),
    @(28,"CC","1051886922","LUIS MIGUEL MERCADO JUNCO","2505",45552,1423500)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Range("B" + $rowNum).Value = $r[1]
    $ws.Range("C" + $rowNum).Value = $r[2]
    $ws.Range("D" + $rowNum).Value = $r[3]
    $ws.Range("E" + $rowNum).Value = $r[4]
    $ws.Range("F" + $rowNum).Value = $r[5]
    $ws.Range("G" + $rowNum).Value = $r[6]
}

# ---------------------------------------------------------------------
# 4) Update the summary header: total Valor Mora, worker/period counts.
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 706056
$ws.Range("C13").Value = 7
$ws.Range("F13").Value = 3
